# Update "想去人数" (F column) values on the "展览" sheet (rows 4-13)
# and on the "全部类型" sheet (rows 5-17) to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 183
$wsExpo.Range("F5").Value = 309
$wsExpo.Range("F6").Value = 400
$wsExpo.Range("F8").Value = 2348
$wsExpo.Range("F9").Value = 390
$wsExpo.Range("F10").Value = 5930
$wsExpo.Range("F11").Value = 148
$wsExpo.Range("F12").Value = 383
$wsExpo.Range("F13").Value = 10

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 183
$wsAll.Range("F6").Value = 309
$wsAll.Range("F7").Value = 400
$wsAll.Range("F11").Value = 2348
$wsAll.Range("F12").Value = 391
$wsAll.Range("F13").Value = 5930
$wsAll.Range("F14").Value = 149
$wsAll.Range("F15").Value = 383
$wsAll.Range("F17").Value = 10
